$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Area"
$ws.Range("AN1").Value = "Comments"

# Match the special header style (s="2") that is already used on D1 (reason_archive)
$ws.Range("D1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
